$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Toggle the "Started" (Yes/No) column for the players whose status flipped.
$ws.Range("C2").Value = "No"
$ws.Range("C4").Value = "No"
$ws.Range("C7").Value = "Yes"
$ws.Range("C9").Value = "Yes"
$ws.Range("C14").Value = "No"
$ws.Range("C15").Value = "No"
$ws.Range("C18").Value = "Yes"
$ws.Range("C20").Value = "No"
$ws.Range("C22").Value = "Yes"
$ws.Range("C24").Value = "Yes"
$ws.Range("C29").Value = "No"
$ws.Range("C31").Value = "No"
$ws.Range("C32").Value = "Yes"
$ws.Range("C34").Value = "Yes"
$ws.Range("C53").Value = "No"
$ws.Range("C55").Value = "No"
$ws.Range("C56").Value = "Yes"
$ws.Range("C58").Value = "Yes"
$ws.Range("C64").Value = "Yes"
$ws.Range("C69").Value = "No"
$ws.Range("C77").Value = "Yes"
$ws.Range("C79").Value = "No"

# Scroll the frozen view back to the top of the data (was parked at row 63
# with C84 selected) so it opens showing the start of the table again.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
